# Update "想去人数" (number of people interested) values for two events
# that are listed on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 200   # was 196
$ws1.Range("F4").Value = 800   # was 798

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 200   # was 196
$ws4.Range("F5").Value = 800   # was 798
